$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03877914339252267
$ws.Range("C2").Value = 0.1223661001864929
$ws.Range("D2").Value = 0.5921565300001895
$ws.Range("E2").Value = 0.2580542788517324
$ws.Range("F2").Value = 0.09150221478225759
$ws.Range("G2").Value = 0.0475100079502194
$ws.Range("H2").Value = 0.1099094607512485
$ws.Range("B3").Value = 0.1894635701346697
$ws.Range("C3").Value = 0.6592539999483663
$ws.Range("D3").Value = 0.3251517487999092
$ws.Range("E3").Value = 0.1585996847304344
$ws.Range("F3").Value = 0.1146074778983962
$ws.Range("G3").Value = 0.1770069306994253
$ws.Range("B4").Value = 0.5408953042689693
$ws.Range("C4").Value = 0.2067930531205122
$ws.Range("D4").Value = 0.04024098905103731
$ws.Range("E4").Value = -0.00375121778100088
$ws.Range("F4").Value = 0.05864823502002824
$ws.Range("G4").Value = 0.02541975869265006
$ws.Range("H4").Value = 0.03155836416756893
$ws.Range("I4").Value = -0.02336658152311553
$ws.Range("J4").Value = -0.02233679963354765
$ws.Range("B5").Value = 0.1798140101391425
$ws.Range("C5").Value = 0.01326194606966766
$ws.Range("D5").Value = -0.03073026076237053
$ws.Range("E5").Value = 0.03166919203865859
$ws.Range("F5").Value = -0.001559284288719586
$ws.Range("G5").Value = 0.004579321186199281
$ws.Range("H5").Value = -0.05034562450448518
$ws.Range("I5").Value = -0.0493158426149173
$ws.Range("B6").Value = 0.2866721711845134
$ws.Range("C6").Value = 0.2426799643524752
$ws.Range("D6").Value = 0.3050794171535043
$ws.Range("E6").Value = 0.2718509408261262
$ws.Range("F6").Value = 0.277989546301045
$ws.Range("G6").Value = 0.2230646006103605
$ws.Range("H6").Value = 0.2240943824999284
$ws.Range("B7").Value = 0.2293351707594228
$ws.Range("C7").Value = 0.2917346235604519
$ws.Range("D7").Value = 0.2585061472330737
$ws.Range("E7").Value = 0.2646447527079926
$ws.Range("F7").Value = 0.2097198070173081
$ws.Range("G7").Value = 0.210749588906876
$ws.Range("B8").Value = 0.03613384424365556
$ws.Range("C8").Value = 0.002905367916277386
$ws.Range("D8").Value = 0.009043973391196253
$ws.Range("E8").Value = -0.04588097229948821
$ws.Range("F8").Value = -0.04485119040992033
$ws.Range("G8").Value = -0.4685329518921513
$ws.Range("H8").Value = 0.02897294484477524
$ws.Range("I8").Value = -0.04257272378961845
$ws.Range("B9").Value = 0.09666542263617448
$ws.Range("C9").Value = 0.1028040281110933
$ws.Range("D9").Value = 0.04787908242040887
$ws.Range("E9").Value = 0.04890886430997676
$ws.Range("F9").Value = -0.3747728971722541
$ws.Range("G9").Value = 0.1227329995646723
$ws.Range("H9").Value = 0.05118733093027864
$ws.Range("B10").Value = 0.05612687574383117
$ws.Range("C10").Value = 0.001201930053146706
$ws.Range("D10").Value = 0.002231711942714591
$ws.Range("E10").Value = -0.4214500495395163
$ws.Range("F10").Value = 0.07605584719741015
$ws.Range("G10").Value = 0.00451017856301647
$ws.Range("B11").Value = -0.06941193049434516
$ws.Range("C11").Value = -0.06838214860477726
$ws.Range("D11").Value = -0.4920639100870082
$ws.Range("E11").Value = 0.005441986649918303
$ws.Range("F11").Value = -0.06610368198447539
$ws.Range("B12").Value = -0.002889654035708851
$ws.Range("C12").Value = -0.4265714155179398
$ws.Range("D12").Value = 0.07093448121898671
$ws.Range("E12").Value = -0.0006111874154069719
$ws.Range("B13").Value = -0.4066743937068965
$ws.Range("C13").Value = 0.09083150303002996
$ws.Range("D13").Value = 0.01928583439563626
$ws.Range("B14").Value = 0.1531996791782531
$ws.Range("C14").Value = 0.08165401054385939
$ws.Range("B15").Value = -0.1180965791298333
